$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) column holds numeric-looking text (e.g. "60.20", "0.999").
# Assigning such a string straight to .Value makes Excel silently reinterpret it
# as a real number (dropping trailing zeros / changing formatting), so a leading
# apostrophe is used to force text entry; Excel stores the apostrophe only as a
# formatting hint, not as part of the cell value. The cell style is then reset to
# "Normal" to drop the quote-prefix styling Excel adds, keeping cell styling as-is.

$ws.Range("D2").Value = '''55.097.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("D3").Value = '''2.314.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''505.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("D6").Value = '''129.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").Value = '''2.308.34'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("D10").Value = '''0.0978'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = '''5.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.53%  '
$ws.Range("D13").Value = '''0.340'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '''23.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.53%  '
$ws.Range("D15").Value = '''2.704.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").Value = '''55.159.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").Value = '''2.281.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").Value = '''4.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").Value = '''312.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '''60.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '''0.995'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  +3.13%  '
$ws.Range("D27").Value = '''7.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("D28").Value = '''172.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").Value = '''0.0₃0712'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.25%  '
$ws.Range("E30").Value = '  +7.01%  '
$ws.Range("E31").Value = '  +3.91%  '
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D34").Value = '''18.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.53%  '
$ws.Range("D35").Value = '''0.996'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  +3.88%  '
$ws.Range("D37").Value = '''0.915'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("D38").Value = '''3.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.25%  '
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("D40").Value = '''1.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.01%  '
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").Value = '''136.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.24%  '
$ws.Range("D43").Value = '''5.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.58%  '
$ws.Range("E44").Value = '  +2.38%  '
$ws.Range("D45").Value = '''260.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.31%  '
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("E47").Value = '  +2.47%  '
$ws.Range("D48").Value = '''0.552'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.56%  '
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("E50").Value = '  +3.23%  '
$ws.Range("D51").Value = '''16.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.35%  '
